$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "addCustomer"
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 15.6666666667
$ws.Columns.Item(3).ColumnWidth = 12.8333333333
$ws.Columns.Item(4).ColumnWidth = 41.6666666667
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"
$ws.Range("A2").Value = "Raman"
$ws.Range("B2").Value = "Arora"
$ws.Range("C2").Value = "A234wd"
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully with customer id"
[void]$ws.Range("A1").Select()
